$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the existing map data from rows 3-8 up to rows 2-7 ---
# Row 2 in the original sheet is empty, so deleting it slides rows 3-8
# up by one (landing exactly on rows 2-7) without disturbing row 1.
$ws.Rows.Item(2).Delete()

# --- Rebuild row 1 as the column-index header (was the old "Syntax..." note) ---
$ws.Range("A1").Clear()
$ws.Range("C1").Clear()

$headerVals = @(0,1,2,3,4,5,6,7,8,9,10,11)
for ($i = 0; $i -lt $headerVals.Length; $i++) {
    $col = 2 + $i   # starts at column B
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headerVals[$i]
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
}
$ws.Rows.Item(1).RowHeight = 19.5

# --- Add the row-index column (column A) for rows 2-14 ---
for ($r = 2; $r -le 14; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $r - 2
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
    $ws.Rows.Item($r).RowHeight = 45
}

# --- Column A is narrower than the rest ---
$ws.Columns.Item(1).ColumnWidth = 3.1666666666666665

# --- Selection / dimension bookkeeping ---
$ws.Range("E9").Select() | Out-Null
